$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add player name "Eric" in column B
$ws.Range("B3").Value = "Eric"

# Row 4: add Points value and update Heure
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "10:15"

# Row 5: update Points value and Heure
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = "10:15"

# Row 6: update Points value and Heure
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "10:15"

# Row 7: update Points value and Heure
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "10:15"

# Row 8: update Points value and Heure
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "10:15"

# Row 9: update Heure
$ws.Range("F9").Value = "10:14"

# Row 10: update Heure
$ws.Range("F10").Value = "10:14"

# Row 11: update Heure
$ws.Range("F11").Value = "10:14"

# Row 12: update Heure
$ws.Range("F12").Value = "10:14"
